$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Switch example group names/numbers for rows 12-16
$ws.Range("B12").Value = "Alice"
$ws.Range("B13").Value = "Bob"
$ws.Range("B14").Value = "Claire"
$ws.Range("B15").Value = "David"
$ws.Range("B16").Value = "Elaine"

# Remove the now-obsolete row 17 entirely
$ws.Rows("17").Delete()

# Update the selection to match the new active cell
$ws.Range("E14").Select()
